$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: DeactivationLeaveBalance scenario
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "39"
$ws.Range("B39").Value = "DeactivationLeaveBalance"
$ws.Range("C39").Value = "DeactivationLeaveBalance"
$ws.Range("D39").Value = "com.darwinbox.leaves.Accural.Custom.DeactivationBalance"
$ws.Range("E39").Value = "deactivation//Leave_Scenarios_Without_Creation.xlsx"
$ws.Range("F39").Value = "All_without_Creation"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "All"

# Row 40: CarryForwardBalance - Custom Leave Cycle scenario
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "40"
$ws.Range("B40").Value = "CarryForwardBalance"
$ws.Range("C40").Value = "Carry Forward Balance -Custom Leave Cycle"
$ws.Range("D40").Value = "com.darwinbox.leaves.Accural.Custom.CarryForwardBalance"
$ws.Range("E40").Value = "Accural//CarryForward.xlsx"
$ws.Range("F40").Value = "All_Scenarios"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "All"

# Row 41: Tenure leave balance scenario
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "41"
$ws.Range("B41").Value = "Tenure"
$ws.Range("C41").Value = "Tenure"
$ws.Range("D41").Value = "com.darwinbox.leaves.TenureLeave.Verify_Tenure_Leave_Balance_of_an_employee"
$ws.Range("E41").Value = "TenureLeaveBalance//Tenure_Leave_Scenarios.xlsx"
$ws.Range("F41").Value = "All_Without_Probation"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "All"

# Row 42: LeaveBalance - 48 employee creation scenario (ClassName entered before TCID)
$ws.Range("D42").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalance_48EmployeeCreation"
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "42"
$ws.Range("B42").Value = "LeaveBalance"
$ws.Range("C42").Value = "LeaveBalance"
$ws.Range("E42").Value = "Accural//LeaveBalance.xlsx"
$ws.Range("F42").Value = "LeaveBalance"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "All"

# Update the view to match the post-edit state: selection moves to A42
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A42").Select()
